# nová logika zoomování, funkční malování
$wb = $excel.ActiveWorkbook

# --- disk_list sheet: add a new row 7 ---
$wsDisk = $wb.Worksheets.Item("disk_list")
$wsDisk.Range("A7").Value = "xfdx"
$wsDisk.Range("B7").Value = "P"
$wsDisk.Range("C7").Value = "\\192.168.000.000\"
$wsDisk.Range("D7").Value = "ss"

# --- Settings sheet: swap B3/B4 values ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B3").Value = 1
$wsSettings.Range("B4").Value = 0
